$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mirror column H (Learning6) values into column I (Learning7) for rows 2-7,
# except row 5 (optimizer) which becomes "nag" instead of "adam".
$ws.Range("I2").Value2 = $ws.Range("H2").Value2
$ws.Range("I3").Value2 = $ws.Range("H3").Value2
$ws.Range("I4").Value2 = $ws.Range("H4").Value2
$ws.Range("I5").Value2 = "nag"
$ws.Range("I6").Value2 = $ws.Range("H6").Value2
$ws.Range("I7").Value2 = $ws.Range("H7").Value2

# Update the active selection on Sheet1 to K10 (from I12).
$ws.Activate()
$ws.Range("K10").Select()
